# ERStudio.xlsx update:
#  - Tsmain path for row 1 (S.No.=1) now points at the TeamServer bootstrapper
#    config file instead of the old abc.txt placeholder.
#  - Column B is widened (best-fit) so the longer path is fully visible.
#  - Selection left on B2 (the cell that was just edited).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Tsmain" value used by row 2 (S.No. 1).
$ws.Range("B2").Value = "\team_server_installer\TeamServer_Bootstrapper\Config.wxi"

# Resize column B to fit the new (much longer) string, and leave the
# selection/active cell on it, matching the saved view state.
$ws.Columns("B:B").ColumnWidth = 52
$ws.Range("B2").Select() | Out-Null
